$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1071.9286
$ws.Range("J17").Value = 1071.9286
$ws.Range("L17").Value = 3215.7858
$ws.Range("N17").Value = -3551.7858
$ws.Range("H33").Value = 49088.617
$ws.Range("I33").Value = 62928.688
$ws.Range("J33").Value = 4800.4
$ws.Range("K33").Value = 62928.688
$ws.Range("L33").Value = 4800.4
$ws.Range("M33").Value = -62699.688
$ws.Range("N33").Value = -5258.4
$ws.Range("H76").Value = 4249.067
$ws.Range("I76").Value = 3759.5557
$ws.Range("J76").Value = 4983.3335
$ws.Range("K76").Value = 3759.5557
$ws.Range("L76").Value = 4983.3335
$ws.Range("M76").Value = -3444.5557
$ws.Range("N76").Value = -5613.3335
$ws.Range("H79").Value = 4249.067
$ws.Range("I79").Value = 3759.5557
$ws.Range("J79").Value = 4983.3335
$ws.Range("K79").Value = 3759.5557
$ws.Range("L79").Value = 4983.3335
$ws.Range("M79").Value = -2667.5557
$ws.Range("N79").Value = -7167.3335
$ws.Range("H132").Value = 11373189
$ws.Range("I132").Value = 14717963
$ws.Range("J132").Value = 959.6
$ws.Range("K132").Value = 44153889
$ws.Range("L132").Value = 2878.8
$ws.Range("M132").Value = -44151359
$ws.Range("N132").Value = -7938.8
$ws.Range("H137").Value = 1035.8524
$ws.Range("I137").Value = 995.05554
$ws.Range("K137").Value = 2985.16662
$ws.Range("M137").Value = -435.16662
$ws.Range("H138").Value = 1814.0303
$ws.Range("I138").Value = 1439.8518
$ws.Range("J138").Value = 3497.8333
$ws.Range("K138").Value = 4319.555399999999
$ws.Range("L138").Value = 10493.4999
$ws.Range("M138").Value = 820.4446000000007
$ws.Range("N138").Value = -20773.4999
$ws.Range("H141").Value = 1246.9518
$ws.Range("I141").Value = 1066.1948
$ws.Range("J141").Value = 3566.6667
$ws.Range("K141").Value = 3198.5844
$ws.Range("L141").Value = 10700.0001
$ws.Range("M141").Value = 1981.4156
$ws.Range("N141").Value = -21060.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 10420.6
$ws.Range("I6").Value = 27302
$ws.Range("J6").Value = 6200.25
$ws.Range("K6").Value = 27302
$ws.Range("L6").Value = 6200.25
$ws.Range("M6").Value = -27129
$ws.Range("N6").Value = -6546.25
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 2000
$ws.Range("K16").Value = 2000
$ws.Range("M16").Value = -1713
$ws.Range("H32").Value = 1413.7
$ws.Range("I32").Value = 1426.0625
$ws.Range("J32").Value = 1117
$ws.Range("K32").Value = 1426.0625
$ws.Range("L32").Value = 1117
$ws.Range("M32").Value = -1139.0625
$ws.Range("N32").Value = -1691
$ws.Range("H37").Value = 12466.667
$ws.Range("J37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15546
$ws.Range("H44").Value = 12780
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 12780
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 12780
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -13756
$ws.Range("H55").Value = 11624
$ws.Range("J55").Value = 12780
$ws.Range("L55").Value = 12780
$ws.Range("N55").Value = -13410
$ws.Range("H61").Value = 1187.7858
$ws.Range("I61").Value = 1124.9429
$ws.Range("K61").Value = 1124.9429
$ws.Range("M61").Value = -912.9429
$ws.Range("H63").Value = 2502.5
$ws.Range("I63").Value = 1005
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 1005
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -319
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 2502.5
$ws.Range("I66").Value = 1005
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 5025
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -1593
$ws.Range("N66").Value = -26864
$ws.Range("H74").Value = 492.77274
$ws.Range("I74").Value = 473.38095
$ws.Range("J74").Value = 900
$ws.Range("K74").Value = 473.38095
$ws.Range("L74").Value = 900
$ws.Range("M74").Value = 400.61905
$ws.Range("N74").Value = -2648
$ws.Range("H77").Value = 492.77274
$ws.Range("I77").Value = 473.38095
$ws.Range("J77").Value = 900
$ws.Range("K77").Value = 2366.90475
$ws.Range("L77").Value = 4500
$ws.Range("M77").Value = 2001.09525
$ws.Range("N77").Value = -13236
$ws.Range("H80").Value = 27489.2
$ws.Range("J80").Value = 27489.2
$ws.Range("L80").Value = 27489.2
$ws.Range("N80").Value = -29485.2
$ws.Range("H83").Value = 27489.2
$ws.Range("J83").Value = 27489.2
$ws.Range("L83").Value = 82467.60000000001
$ws.Range("N83").Value = -92451.60000000001
$ws.Range("H136").Value = 1187.7858
$ws.Range("I136").Value = 1124.9429
$ws.Range("K136").Value = 3374.8287
$ws.Range("M136").Value = -824.8287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18175.084
$ws.Range("J82").Value = 28972.715
$ws.Range("L82").Value = 28972.715
$ws.Range("N82").Value = -29738.715
$ws.Range("H85").Value = 18175.084
$ws.Range("J85").Value = 28972.715
$ws.Range("L85").Value = 28972.715
$ws.Range("N85").Value = -31624.715
$ws.Range("H86").Value = 72769.44
$ws.Range("I86").Value = 88747.08
$ws.Range("J86").Value = 3533
$ws.Range("K86").Value = 88747.08
$ws.Range("L86").Value = 3533
$ws.Range("M86").Value = -87624.08
$ws.Range("N86").Value = -5779
$ws.Range("H89").Value = 72769.44
$ws.Range("I89").Value = 88747.08
$ws.Range("J89").Value = 3533
$ws.Range("K89").Value = 443735.4
$ws.Range("L89").Value = 17665
$ws.Range("M89").Value = -438119.4
$ws.Range("N89").Value = -28897
$ws.Range("H94").Value = 466.63635
$ws.Range("I94").Value = 413.3
$ws.Range("K94").Value = 413.3
$ws.Range("M94").Value = 37.69999999999999
$ws.Range("H107").Value = 100046690
$ws.Range("I107").Value = 200091920
$ws.Range("K107").Value = 200091920
$ws.Range("M107").Value = -200090000
$ws.Range("H134").Value = 2230.647
$ws.Range("I134").Value = 1975.26
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 5925.78
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -3390.78
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 2750
$ws.Range("I23").Value = 2500
$ws.Range("J23").Value = 3000
$ws.Range("K23").Value = 2500
$ws.Range("L23").Value = 3000
$ws.Range("M23").Value = -2260
$ws.Range("N23").Value = -3480
$ws.Range("H27").Value = 2750
$ws.Range("I27").Value = 2500
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 2500
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -2308
$ws.Range("N27").Value = -3384
$ws.Range("H134").Value = 1302.2142
$ws.Range("I134").Value = 1100.6842
$ws.Range("K134").Value = 3302.0526
$ws.Range("M134").Value = -767.0526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8060.2764
$ws.Range("J131").Value = 9043.313
$ws.Range("L131").Value = 27129.939
$ws.Range("N131").Value = -37209.939

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 32500
$ws.Range("J32").Value = 32500
$ws.Range("L32").Value = 32500
$ws.Range("N32").Value = -33092
$ws.Range("H43").Value = 3361.3845
$ws.Range("I43").Value = 808.5714
$ws.Range("J43").Value = 6339.6665
$ws.Range("K43").Value = 808.5714
$ws.Range("L43").Value = 6339.6665
$ws.Range("M43").Value = -657.5714
$ws.Range("N43").Value = -6641.6665
$ws.Range("H46").Value = 10899.3
$ws.Range("I46").Value = 8000
$ws.Range("J46").Value = 11221.444
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 11221.444
$ws.Range("M46").Value = -7844
$ws.Range("N46").Value = -11533.444
$ws.Range("H110").Value = 38900
$ws.Range("J110").Value = 38900
$ws.Range("L110").Value = 38900
$ws.Range("N110").Value = -47080
$ws.Range("H113").Value = 1683.3334
$ws.Range("I113").Value = 1858.25
$ws.Range("J113").Value = 1575.6923
$ws.Range("K113").Value = 1858.25
$ws.Range("L113").Value = 1575.6923
$ws.Range("M113").Value = 311.75
$ws.Range("N113").Value = -5915.6923
$ws.Range("H122").Value = 614.2632
$ws.Range("I122").Value = 527.7059
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 1583.1177
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = 866.8822999999998
$ws.Range("N122").Value = -8950
$ws.Range("H126").Value = 2880.9412
$ws.Range("I126").Value = 3045.6667
$ws.Range("J126").Value = 2485.6
$ws.Range("K126").Value = 9137.000100000001
$ws.Range("L126").Value = 7456.799999999999
$ws.Range("M126").Value = -6667.000100000001
$ws.Range("N126").Value = -12396.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2546
$ws.Range("I100").Value = 1862.5
$ws.Range("J100").Value = 5280
$ws.Range("K100").Value = 1862.5
$ws.Range("L100").Value = 5280
$ws.Range("M100").Value = -1321.5
$ws.Range("N100").Value = -6362
$ws.Range("H136").Value = 1140.3265
$ws.Range("I136").Value = 941.32556
$ws.Range("K136").Value = 2823.97668
$ws.Range("M136").Value = -273.9766799999998
